# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    4  = 2179
    5  = 68
    6  = 12759
    8  = 109
    9  = 505
    10 = 467
    12 = 958
    13 = 13683
    14 = 14086
    19 = 22
    23 = 1064
    24 = 109
    26 = 929
    27 = 5201
    29 = 266
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
